# Apply the eDNA ASV reassignment edits.
# Swap A:D text values between rows 42/43 and between rows 55/56,
# and move the J-column value (0) from row 42 to row 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 42 <-> 43 ---
$ws.Range("A42").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B42").Value = "Homo sapiens"
$ws.Range("C42").Value = "Human"
$ws.Range("D42").Value = "Human"

$ws.Range("A43").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B43").Value = "unassigned"
$ws.Range("C43").Value = "unassigned"
$ws.Range("D43").Value = "unassigned"

# J42 had a value of 0; it moves to J43 and J42 becomes empty
$ws.Range("J42").ClearContents()
$ws.Range("J43").Value = 0

# --- Rows 55 <-> 56 ---
$ws.Range("A55").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B55").Value = "unassigned"
$ws.Range("C55").Value = "unassigned"
$ws.Range("D55").Value = "unassigned"

$ws.Range("A56").Value = "975b1dbdc7405f6e27bf63893e91e0ed"
$ws.Range("B56").Value = "Centropristis striata"
$ws.Range("C56").Value = "Black sea bass"
$ws.Range("D56").Value = "Teleost Fish"
